$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B..N -> C..N, old F.. -> G..)
$ws.Columns("B:B").Insert()

# The old row-1 placeholder cell (F1) is now G1 after the column insert;
# it carried no content, only a style - clear it out entirely so row 1 disappears.
$ws.Range("G1").Clear()

# New "ano" column header (new shared string) + values for the three data rows
$ws.Range("B2").Value = "ano"
$ws.Range("B3").Value = 2023
$ws.Range("B4").Value = 2023
$ws.Range("B5").Value = 2023

# Style B2 like the other header cells in col A (bold header look is already
# inherited from the column insert - A2's style carried over to B2).

# Give B3:B5 a dedicated style: bold Calibri 9, centered, top-aligned, wrap
# text, solid fill FBC995, no border (matches the other header/value palette
# used throughout the sheet).
$vals = $ws.Range("B3:B5")
$vals.Font.Name = "Calibri"
$vals.Font.Size = 9
$vals.Font.Bold = $true
$vals.Interior.Color = 9816571
$vals.HorizontalAlignment = -4108
$vals.VerticalAlignment = -4160
$vals.WrapText = $true
$vals.Borders.LineStyle = -4142

# Match column B's width to column A's so the grouped <cols> entry lines up
$ws.Range("B1").ColumnWidth = $ws.Range("A1").ColumnWidth

# Leave the cursor where the author ended up after the edits
$null = $ws.Range("B6").Select()
